$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.229.70'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '3.460.44'
$ws.Range('E3').Value = '  -1.37%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.609'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.54%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '3.458.75'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.139'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.94'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.74%  '
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').Value = '4.052.64'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '31.99'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '67.203.16'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000177'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').Value = '3.459.66'
$ws.Range('E18').Value = '  -1.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '392.72'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.538'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '71.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.175'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.11'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.48'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.32'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.58'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.878'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.81'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.87'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.22%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.65'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.64%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0718'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.10%  '
$ws.Range('D46').Value = '2.753.91'
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.29'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.36'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.74%  '
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '325.88'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.72%  '
$ws.Range('E51').Value = '  -3.46%  '
